$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.104.89"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "2.966.79"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'596.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'149.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "2.966.56"
$ws.Range("D9").Value = "'0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.17%  "
$ws.Range("D12").Value = "'0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.03%  "
$ws.Range("D14").Value = "'33.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "3.462.85"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "63.099.05"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "'6.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "2.946.00"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'445.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").Value = "'13.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "'0.674"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'7.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'11.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").Value = "'81.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "'11.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000106"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +20.18%  "
$ws.Range("D31").Value = "'7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.33%  "
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "'26.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").Value = "'0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'3.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.19%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").Value = "'5.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'2.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").Value = "'49.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "'8.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "'41.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "2.716.42"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'370.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").Value = "'134.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D50").Value = "'23.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
